{"js": "// Word Online / Office.js (Word JS API) edit script.\n// Body of: async (context) => { ... }\n//\n// This reproduces the authoring change described in the task's diff:\n//   1. The three chart drawings (wp:inline pictures) get \"noProof\" plus\n//      an explicit en-CA language stamped on their run properties (this is\n//      what Word does when it re-measures/re-proofs inline objects after\n//      an edit elsewhere in the document).\n//   2. A new (unused, hidden-until-used) \"FollowedHyperlink\" character\n//      style is added to the style sheet, mirroring the built-in\n//      \"Hyperlink\" style that was already present.\n\n// --- 1. Stamp noProof + en-CA language on every chart's anchor run -------\nconst pics = context.document.body.inlinePictures;\npics.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < pics.items.length; i++) {\n  const rng = pics.items[i].getRange();\n  rng.hasNoProofing = true;\n  rng.languageId = \"en-CA\";\n  rng.languageIdFarEast = \"en-CA\";\n}\nawait context.sync();\n\n// --- 2. Add the FollowedHyperlink character style ------------------------\ncontext.document.addStyle(\"FollowedHyperlink\", Word.StyleType.character);\nawait context.sync();\n\nconst styles = context.document.getStyles();\nconst followedHyperlink = styles.getByName(\"FollowedHyperlink\");\nfollowedHyperlink.baseStyle = \"DefaultParagraphFont\";\nfollowedHyperlink.priority = 99;\nfollowedHyperlink.unhideWhenUsed = true;\nfollowedHyperlink.font.color = \"#000000\";\nfollowedHyperlink.font.underline = Word.UnderlineType.single;\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d / $doc are pre-seeded by the host.\n#\n# This reproduces the authoring change described in the task's diff:\n#   1. The three chart drawings (InlineShapes) get \"NoProofing\" plus an\n#      explicit en-CA language stamped on their anchor run - this is what\n#      Word does when it re-measures/re-proofs inline objects after an\n#      edit elsewhere in the document.\n#   2. A new (unused, hidden-until-used) \"FollowedHyperlink\" character\n#      style is added to the style sheet, mirroring the built-in\n#      \"Hyperlink\" style that was already present.\n\n$d = $word.ActiveDocument\n\n# --- 1. Stamp NoProofing + en-CA language on every chart's anchor run ----\n$shapeCount = $d.InlineShapes.Count\nfor ($i = 1; $i -le $shapeCount; $i++) {\n  $shp = $d.InlineShapes.Item($i)\n  $rng = $shp.Range\n  $rng.NoProofing = 1\n  $rng.LanguageID = \"en-CA\"\n  $rng.LanguageIDFarEast = \"en-CA\"\n}\n\n# --- 2. Add the FollowedHyperlink character style -------------------------\n$followedHyperlink = $d.Styles.Add(\"FollowedHyperlink\", 2)\n$followedHyperlink.BaseStyle = \"DefaultParagraphFont\"\n$followedHyperlink.Priority = 99\n$followedHyperlink.UnhideWhenUsed = 1\n$followedHyperlink.Font.Color = 0\n$followedHyperlink.Font.Underline = 1\n"}
